$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage so numeric-looking strings (e.g. trailing zeros,
# leading zeros, percent signs) keep their exact original text formatting
# instead of being reinterpreted as numbers.
$cellNames = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","E15","D16","E16","D17","E17","D19","E19","E20","D21","E21","E22","D23","E23","E24","D25","E25","D26","E26","D27","E27","D39","E39","D40","E40","D41","E41","E42","D43","E43","D44","E44","E45","D46","E46","D47","E47","D48","E48","E49","D50","E50","D51","E51")
foreach ($name in $cellNames) {
    $ws.Range($name).NumberFormat = "@"
}

$ws.Range('D2').Value = '329.71'
$ws.Range('E2').Value = '1.29%'
$ws.Range('D3').Value = '41.44'
$ws.Range('E3').Value = '4.43%'
$ws.Range('D4').Value = '5.646'
$ws.Range('E4').Value = '-1.02%'
$ws.Range('D5').Value = '0.08274'
$ws.Range('E5').Value = '2.95%'
$ws.Range('D6').Value = '2.042'
$ws.Range('E6').Value = '0.18%'
$ws.Range('D7').Value = '8.763'
$ws.Range('E7').Value = '1.45%'
$ws.Range('D8').Value = '4.528'
$ws.Range('E8').Value = '0.64%'
$ws.Range('D10').Value = '0.9252'
$ws.Range('E10').Value = '0.10%'
$ws.Range('D11').Value = '0.1275'
$ws.Range('E11').Value = '1.04%'
$ws.Range('D12').Value = '0.1956'
$ws.Range('E12').Value = '-0.04%'
$ws.Range('D13').Value = '0.09421'
$ws.Range('E13').Value = '2.64%'
$ws.Range('D14').Value = '0.03973'
$ws.Range('E14').Value = '11.56%'
$ws.Range('E15').Value = '1.14%'
$ws.Range('D16').Value = '0.001309'
$ws.Range('E16').Value = '1.99%'
$ws.Range('D17').Value = '0.006100'
$ws.Range('E17').Value = '-2.72%'
$ws.Range('D19').Value = '3.445'
$ws.Range('E19').Value = '2.38%'
$ws.Range('E20').Value = '0.01%'
$ws.Range('D21').Value = '8.370'
$ws.Range('E21').Value = '-4.46%'
$ws.Range('E22').Value = '1.81%'
$ws.Range('D23').Value = '0.2662'
$ws.Range('E23').Value = '-0.21%'
$ws.Range('E24').Value = '-0.18%'
$ws.Range('D25').Value = '0.001255'
$ws.Range('E25').Value = '-0.32%'
$ws.Range('D26').Value = '0.004316'
$ws.Range('E26').Value = '-6.28%'
$ws.Range('D27').Value = '0.0001201'
$ws.Range('E27').Value = '0.86%'
$ws.Range('D39').Value = '0.02768'
$ws.Range('E39').Value = '10.75%'
$ws.Range('D40').Value = '0.05510'
$ws.Range('E40').Value = '3.73%'
$ws.Range('D41').Value = '0.007916'
$ws.Range('E41').Value = '5.86%'
$ws.Range('E42').Value = '0.85%'
$ws.Range('D43').Value = '0.008942'
$ws.Range('E43').Value = '-9.72%'
$ws.Range('D44').Value = '0.002142'
$ws.Range('E44').Value = '1.20%'
$ws.Range('E45').Value = '0.24%'
$ws.Range('D46').Value = '0.00007012'
$ws.Range('E46').Value = '4.90%'
$ws.Range('D47').Value = '0.00000000751'
$ws.Range('E47').Value = '0.04%'
$ws.Range('D48').Value = '0.003192'
$ws.Range('E48').Value = '4.95%'
$ws.Range('E49').Value = '0.16%'
$ws.Range('D50').Value = '0.00002102'
$ws.Range('E50').Value = '0.04%'
$ws.Range('D51').Value = '0.0002002'
$ws.Range('E51').Value = '0.04%'
